# Update "Generate Report for Handback" timestamps across the workbook.

$wb = $excel.ActiveWorkbook

# Sheet "Overview": Latest HO Xliff Generate Date (column G, row 2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-13 23:29:47"

# Sheet "zh-cn": Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2), row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-13 23:29:39"
$wsZhCn.Range("K2").Value = "2016-08-13 23:30:18"

# Sheet "de-de": Correspond Handback DateTime (K2), row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-13 23:30:28"
